$d = $word.ActiveDocument

# Locate the "Travis CI" section: it starts with the bold "Travis CI" heading
# paragraph and runs through the trailing "Src - https://en.wikipedia.org/..."
# paragraph. That whole section (and the blank paragraph inside it) is being
# removed, but the "_GoBack" bookmark that lived at the end of the heading
# paragraph needs to survive, anchored in the blank paragraph that precedes
# the section (which becomes the document's final paragraph).

$startIndex = -1
$endIndex = -1

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $text = $d.Paragraphs.Item($i).Range.Text.Trim()
    if ($text -eq "Travis CI") {
        $startIndex = $i
    }
    if ($text -eq "Src - https://en.wikipedia.org/wiki/Travis_CI") {
        $endIndex = $i
    }
}

if ($startIndex -gt 0 -and $endIndex -ge $startIndex) {
    $keepIndex = $startIndex - 1

    $sectionStart = $d.Paragraphs.Item($startIndex).Range.Start
    $sectionEnd = $d.Paragraphs.Item($endIndex).Range.End

    $sectionRange = $d.Range($sectionStart, $sectionEnd)
    $sectionRange.Delete()

    $keepRange = $d.Paragraphs.Item($keepIndex).Range
    $d.Bookmarks.Add("_GoBack", $keepRange)
}
